# Add 8 new travel rows (563-570) to the "travels" worksheet, matching the
# diff that extends the sheet's data range from A1:D562 to A1:D570.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("travels")

$newRows = @(
    @("18_Car-CA",     "Cairo", "Giza",     6.9),
    @("18_Car-GIZ",    "Giza",  "Qalyubia", 31.8),
    @("180_Car-CA",    "Cairo", "Giza",     6.9),
    @("180_Car-GIZ",   "Giza",  "Qalyubia", 31.8),
    @("181_Car-CA",    "Cairo", "Giza",     6.9),
    @("181_Car-GIZ",   "Giza",  "Qalyubia", 31.8),
    @("123456_Car-CA", "Cairo", "Giza",     6.9),
    @("123456_Car-GIZ","Giza",  "Qalyubia", 31.8)
)

$startRow = 563
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
    $ws.Cells.Item($r, 4).Value = $newRows[$i][3]
}
